# edit.ps1 - applies the Fukushima.docx content edits described by the diff.
#
# Summary of changes:
#  1. Split the italic run "increased the Earth's rotational speed by "
#     into "increased the" + "_GoBack" bookmark + " Earth's rotational speed by "
#     (this also moves the pre-existing "_GoBack" bookmark away from its old
#     location near "habitat destruction (", since Word only keeps one).
#  2. Remove the "Nuclear is unique. " sentence/run and change
#     "nuclear energy is a source" -> "nuclear is a source" in the next run.
#  3. Split the "...sprawling solar farms? Perhaps it " run, removing the
#     trailing "Perhaps it " text and inserting a new sentence about hydro
#     energy & dams, then re-appending " Perhaps it " as its own run.
#  4. Insert the word " proper" before " public discussion".
#  5. (handled automatically by #1) old "_GoBack" bookmark removed.

$d = $word.ActiveDocument

function New-TempSplit {
    param($Position)
    # Adding then immediately deleting a bookmark at a collapsed range forces
    # the run that spans that position to be split into two runs at that
    # exact character offset, without altering any visible content.
    $collapsed = $d.Range($Position, $Position)
    $d.Bookmarks.Add("__tmp_split__", $collapsed) | Out-Null
    $d.Bookmarks("__tmp_split__").Delete()
}

# ---------------------------------------------------------------------
# Change 1: "increased the Earth's rotational speed by " ->
#           "increased the" + _GoBack bookmark + " Earth's rotational speed by "
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("increased the")
$splitPos = $r1.End
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos)) | Out-Null

# ---------------------------------------------------------------------
# Change 2: remove "Nuclear is unique. " and tweak following sentence.
# ---------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Nuclear is unique. ")
$r2.Text = ""

$r2b = $d.Content
$r2b.Find.Execute("Unlike fossil fuels such as gas and coal, nuclear energy is a source")
$r2b.Text = "Unlike fossil fuels such as gas and coal, nuclear is a source"

# ---------------------------------------------------------------------
# Change 3: split off "Perhaps it " and insert the hydro/dams sentence
#           between the solar-farms question and "Perhaps it ".
# ---------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("sprawling solar farms? Perhaps it ")
$r3.Text = "sprawling solar farms?"

$r3b = $d.Content
$r3b.Find.Execute("sprawling solar farms?")
$insPos = $r3b.End
New-TempSplit -Position $insPos
$ins = $d.Range($insPos, $insPos)
$ins.InsertAfter(" Is it right to destroy natural water courses with dams for hydro energy?")

$r3c = $d.Content
$r3c.Find.Execute("hydro energy?")
$insPos2 = $r3c.End
New-TempSplit -Position $insPos2
$ins2 = $d.Range($insPos2, $insPos2)
$ins2.InsertAfter(" Perhaps it ")

# ---------------------------------------------------------------------
# Change 4: insert " proper" before " public discussion".
# ---------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("to support this idea, nor has there been the public discussion")
$r4.Text = "to support this idea, nor has there been the"

$r4b = $d.Content
$r4b.Find.Execute("to support this idea, nor has there been the")
$insPos3 = $r4b.End
New-TempSplit -Position $insPos3
$ins3 = $d.Range($insPos3, $insPos3)
$ins3.InsertAfter(" proper")

$r4c = $d.Content
$r4c.Find.Execute(" proper")
$insPos4 = $r4c.End
New-TempSplit -Position $insPos4
$ins4 = $d.Range($insPos4, $insPos4)
$ins4.InsertAfter(" public discussion")

Write-Host "Done applying Fukushima.docx edits."
